$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.417, -0.367, 0.237, -85, -5, -94),
    @(-0.177, -0.276, 0.441, -90, 0, 56),
    @(0.045, 0.666, 0.146, -89, 4, -77),
    @(0.437, -0.575, 0.058, -34, 20, 15),
    @(-0.305, -0.08699999999999999, 0.716, 51, 31, -21),
    @(-0.499, -0.034, -0.131, 40, -61, 44),
    @(0.028, -0.255, 0.522, 53, -30, 175),
    @(-0.378, -0.29, -0.19, 71, 5, 130),
    @(-0.04, 0.334, 0.605, 104, -31, -126),
    @(0.07199999999999999, -0.3, 0.579, -165, 11, -144)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$wb.Save()
